# Daily attendance processing - 2025-10-12 08:45:14
# For every cell in the "Recorded By" column (G) whose value is a
# comma-separated list starting with "System", reverse the order of the
# list (e.g. "System, foo@bar.com" -> "foo@bar.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -and $text.StartsWith("System,")) {
        $parts = $text.Split(",")
        $count = $parts.Length

        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i].Trim()
        }

        $newValue = [string]::Join(", ", $reversed)
        $cell.Value = $newValue
    }
}
